# Penalty Reward System (unfinished) - remove some historical rows from the
# weekly/monthly PO data sheets. Rows are deleted from the bottom up so that
# row indices of not-yet-processed rows remain valid.

$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": delete rows 7, 4, 3, 2 (bottom to top) ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows(7).Delete()
$ws1.Rows(4).Delete()
$ws1.Rows(3).Delete()
$ws1.Rows(2).Delete()

# --- Sheet "Monthly Trend": delete rows 4, 2 (bottom to top) ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Rows(4).Delete()
$ws2.Rows(2).Delete()
